$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): shift label values B<-C<-D, with D getting "recall"
$ws.Range("B1").Value = "accuracy"
$ws.Range("C1").Value = "precision"
$ws.Range("D1").Value = "recall"

# Row 2: swap B2 and C2 values (D2 stays the same)
$ws.Range("B2").Value = 91.66666666666666
$ws.Range("C2").Value = 84.02777777777779

# Row 3: swap B3 and C3 values (D3 stays the same)
$ws.Range("B3").Value = 91.66666666666666
$ws.Range("C3").Value = 84.02777777777779
